# coa_db sheet: the record for "SKU 2" (row 3) is being dropped from the
# template, and the status of the remaining record ("SKU 1", row 2) moves
# from "Viejo" to "Nuevo" while the values are positioned for the PDF
# template (per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status column (I) for the "SKU 1" row: Viejo -> Nuevo
$ws.Range("I2").Value = "Nuevo"

# Remove the second data row ("SKU 2" / lote 0002) entirely.
$ws.Rows(3).Delete() | Out-Null

# Leave the cursor where the author was last working.
$ws.Range("E12").Select() | Out-Null
